# Applies updated profit-calculation values to the per-class Leve profit tables
# (sheets ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) produced by the scheduled
# market-data refresh. Each table has columns:
#   H=currentAveragePrice I=currentAveragePriceNQ J=currentAveragePriceHQ
#   K=LevePriceNQ L=LevePriceHQ M=LeveProfitNQ N=LeveProfitHQ
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(17, 8).Value = 1138.7407  # H17
$ws.Cells.Item(17, 10).Value = 1364.409  # J17
$ws.Cells.Item(17, 12).Value = 4093.227  # L17
$ws.Cells.Item(17, 14).Value = -4429.227000000001  # N17
$ws.Cells.Item(18, 8).Value = 7064.3335  # H18
$ws.Cells.Item(18, 9).Value = 10024.5  # I18
$ws.Cells.Item(18, 10).Value = 1144  # J18
$ws.Cells.Item(18, 11).Value = 10024.5  # K18
$ws.Cells.Item(18, 12).Value = 1144  # L18
$ws.Cells.Item(18, 13).Value = -9740.5  # M18
$ws.Cells.Item(18, 14).Value = -1712  # N18
$ws.Cells.Item(39, 8).Value = 2198.7273  # H39
$ws.Cells.Item(39, 9).Value = 1375.3334  # I39
$ws.Cells.Item(39, 10).Value = 5904  # J39
$ws.Cells.Item(39, 11).Value = 4126.0002  # K39
$ws.Cells.Item(39, 12).Value = 17712  # L39
$ws.Cells.Item(39, 13).Value = -3830.0002  # M39
$ws.Cells.Item(39, 14).Value = -18304  # N39
$ws.Cells.Item(43, 8).Value = 8207.691999999999  # H43
$ws.Cells.Item(43, 10).Value = 3100  # J43
$ws.Cells.Item(43, 12).Value = 3100  # L43
$ws.Cells.Item(43, 14).Value = -3238  # N43
$ws.Cells.Item(52, 8).Value = 298.94116  # H52
$ws.Cells.Item(52, 9).Value = 298  # I52
$ws.Cells.Item(52, 11).Value = 894  # K52
$ws.Cells.Item(52, 13).Value = -734  # M52
$ws.Cells.Item(58, 8).Value = 3749.375  # H58
$ws.Cells.Item(58, 9).Value = 1199  # I58
$ws.Cells.Item(58, 10).Value = 8000  # J58
$ws.Cells.Item(58, 11).Value = 3597  # K58
$ws.Cells.Item(58, 12).Value = 24000  # L58
$ws.Cells.Item(58, 13).Value = -3447  # M58
$ws.Cells.Item(58, 14).Value = -24300  # N58
$ws.Cells.Item(62, 8).Value = 111142296  # H62
$ws.Cells.Item(62, 9).Value = 250001740  # I62
$ws.Cells.Item(62, 10).Value = 54731.2  # J62
$ws.Cells.Item(62, 11).Value = 250001740  # K62
$ws.Cells.Item(62, 12).Value = 54731.2  # L62
$ws.Cells.Item(62, 13).Value = -250001116  # M62
$ws.Cells.Item(62, 14).Value = -55979.2  # N62
$ws.Cells.Item(65, 8).Value = 111142296  # H65
$ws.Cells.Item(65, 9).Value = 250001740  # I65
$ws.Cells.Item(65, 10).Value = 54731.2  # J65
$ws.Cells.Item(65, 11).Value = 1250008700  # K65
$ws.Cells.Item(65, 12).Value = 273656  # L65
$ws.Cells.Item(65, 13).Value = -1250005580  # M65
$ws.Cells.Item(65, 14).Value = -279896  # N65
$ws.Cells.Item(86, 8).Value = 99538370  # H86
$ws.Cells.Item(86, 9).Value = 88890296  # I86
$ws.Cells.Item(86, 10).Value = 117285160  # J86
$ws.Cells.Item(86, 11).Value = 88890296  # K86
$ws.Cells.Item(86, 12).Value = 117285160  # L86
$ws.Cells.Item(86, 13).Value = -88889173  # M86
$ws.Cells.Item(86, 14).Value = -117287406  # N86
$ws.Cells.Item(89, 8).Value = 99538370  # H89
$ws.Cells.Item(89, 9).Value = 88890296  # I89
$ws.Cells.Item(89, 10).Value = 117285160  # J89
$ws.Cells.Item(89, 11).Value = 444451480  # K89
$ws.Cells.Item(89, 12).Value = 586425800  # L89
$ws.Cells.Item(89, 13).Value = -444445864  # M89
$ws.Cells.Item(89, 14).Value = -586437032  # N89
$ws.Cells.Item(121, 8).Value = 5314.2856  # H121
$ws.Cells.Item(121, 10).Value = 5314.2856  # J121
$ws.Cells.Item(121, 12).Value = 15942.8568  # L121
$ws.Cells.Item(121, 14).Value = -19436.8568  # N121
$ws.Cells.Item(125, 8).Value = 20833880  # H125
$ws.Cells.Item(125, 9).Value = 31250498  # I125
$ws.Cells.Item(125, 10).Value = 643  # J125
$ws.Cells.Item(125, 11).Value = 281254482  # K125
$ws.Cells.Item(125, 12).Value = 5787  # L125
$ws.Cells.Item(125, 13).Value = -281252022  # M125
$ws.Cells.Item(125, 14).Value = -10707  # N125
$ws.Cells.Item(137, 8).Value = 2357.2  # H137
$ws.Cells.Item(137, 9).Value = 2339.1304  # I137
$ws.Cells.Item(137, 10).Value = 2372.5925  # J137
$ws.Cells.Item(137, 11).Value = 7017.3912  # K137
$ws.Cells.Item(137, 12).Value = 7117.7775  # L137
$ws.Cells.Item(137, 13).Value = -4467.3912  # M137
$ws.Cells.Item(137, 14).Value = -12217.7775  # N137
$ws.Cells.Item(138, 8).Value = 1731064.6  # H138
$ws.Cells.Item(138, 9).Value = 4549.25  # I138
$ws.Cells.Item(138, 10).Value = 2007307.1  # J138
$ws.Cells.Item(138, 11).Value = 13647.75  # K138
$ws.Cells.Item(138, 12).Value = 6021921.300000001  # L138
$ws.Cells.Item(138, 13).Value = -8507.75  # M138
$ws.Cells.Item(138, 14).Value = -6032201.300000001  # N138
$ws.Cells.Item(141, 8).Value = 14999  # H141
$ws.Cells.Item(141, 9).Value = 0  # I141
$ws.Cells.Item(141, 11).Value = 0  # K141
$ws.Cells.Item(141, 13).ClearContents()  # M141

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 3598.8608  # H32
$ws.Cells.Item(32, 9).Value = 3021.861  # I32
$ws.Cells.Item(32, 10).Value = 9533.714  # J32
$ws.Cells.Item(32, 11).Value = 3021.861  # K32
$ws.Cells.Item(32, 12).Value = 9533.714  # L32
$ws.Cells.Item(32, 13).Value = -2734.861  # M32
$ws.Cells.Item(32, 14).Value = -10107.714  # N32
$ws.Cells.Item(51, 8).Value = 88571.42999999999  # H51
$ws.Cells.Item(51, 9).Value = 80000  # I51
$ws.Cells.Item(51, 10).Value = 110000  # J51
$ws.Cells.Item(51, 11).Value = 80000  # K51
$ws.Cells.Item(51, 12).Value = 110000  # L51
$ws.Cells.Item(51, 13).Value = -79244  # M51
$ws.Cells.Item(51, 14).Value = -111512  # N51
$ws.Cells.Item(53, 8).Value = 0  # H53
$ws.Cells.Item(53, 9).Value = 0  # I53
$ws.Cells.Item(53, 11).Value = 0  # K53
$ws.Cells.Item(53, 13).ClearContents()  # M53
$ws.Cells.Item(61, 8).Value = 6371.8184  # H61
$ws.Cells.Item(61, 9).Value = 3297.0435  # I61
$ws.Cells.Item(61, 11).Value = 3297.0435  # K61
$ws.Cells.Item(61, 13).Value = -3085.0435  # M61
$ws.Cells.Item(74, 8).Value = 66921.66  # H74
$ws.Cells.Item(74, 9).Value = 125296.54  # I74
$ws.Cells.Item(74, 11).Value = 125296.54  # K74
$ws.Cells.Item(74, 13).Value = -124422.54  # M74
$ws.Cells.Item(75, 8).Value = 0  # H75
$ws.Cells.Item(75, 10).Value = 0  # J75
$ws.Cells.Item(75, 12).ClearContents()  # L75
$ws.Cells.Item(75, 14).ClearContents()  # N75
$ws.Cells.Item(77, 8).Value = 66921.66  # H77
$ws.Cells.Item(77, 9).Value = 125296.54  # I77
$ws.Cells.Item(77, 11).Value = 626482.7  # K77
$ws.Cells.Item(77, 13).Value = -622114.7  # M77
$ws.Cells.Item(78, 8).Value = 0  # H78
$ws.Cells.Item(78, 10).Value = 0  # J78
$ws.Cells.Item(78, 12).ClearContents()  # L78
$ws.Cells.Item(78, 14).ClearContents()  # N78
$ws.Cells.Item(93, 8).Value = 0  # H93
$ws.Cells.Item(93, 10).Value = 0  # J93
$ws.Cells.Item(93, 12).ClearContents()  # L93
$ws.Cells.Item(93, 14).ClearContents()  # N93
$ws.Cells.Item(124, 8).Value = 52501.668  # H124
$ws.Cells.Item(124, 10).Value = 52501.668  # J124
$ws.Cells.Item(124, 12).Value = 52501.668  # L124
$ws.Cells.Item(124, 14).Value = -62321.668  # N124
$ws.Cells.Item(132, 8).Value = 1320440.2  # H132
$ws.Cells.Item(132, 9).Value = 2026121.2  # I132
$ws.Cells.Item(132, 10).Value = 9889.857  # J132
$ws.Cells.Item(132, 11).Value = 6078363.6  # K132
$ws.Cells.Item(132, 12).Value = 29669.571  # L132
$ws.Cells.Item(132, 13).Value = -6075833.6  # M132
$ws.Cells.Item(132, 14).Value = -34729.571  # N132
$ws.Cells.Item(135, 8).Value = 97142  # H135
$ws.Cells.Item(135, 10).Value = 97142  # J135
$ws.Cells.Item(135, 12).Value = 97142  # L135
$ws.Cells.Item(135, 14).Value = -107282  # N135
$ws.Cells.Item(136, 8).Value = 6371.8184  # H136
$ws.Cells.Item(136, 9).Value = 3297.0435  # I136
$ws.Cells.Item(136, 11).Value = 9891.130500000001  # K136
$ws.Cells.Item(136, 13).Value = -7341.130500000001  # M136
$ws.Cells.Item(139, 8).Value = 90000  # H139
$ws.Cells.Item(139, 10).Value = 90000  # J139
$ws.Cells.Item(139, 12).Value = 90000  # L139
$ws.Cells.Item(139, 14).Value = -100280  # N139

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 5378057.5  # H20
$ws.Cells.Item(20, 10).Value = 1887.2858  # J20
$ws.Cells.Item(20, 12).Value = 1887.2858  # L20
$ws.Cells.Item(20, 14).Value = -2381.2858  # N20
$ws.Cells.Item(25, 8).Value = 4020  # H25
$ws.Cells.Item(25, 9).Value = 0  # I25
$ws.Cells.Item(25, 10).Value = 4020  # J25
$ws.Cells.Item(25, 11).Value = 0  # K25
$ws.Cells.Item(25, 12).ClearContents()  # L25
$ws.Cells.Item(25, 13).ClearContents()  # M25
$ws.Cells.Item(25, 14).Value = -4490  # N25
$ws.Cells.Item(28, 8).Value = 48738.5  # H28
$ws.Cells.Item(28, 10).Value = 48738.5  # J28
$ws.Cells.Item(28, 12).Value = 48738.5  # L28
$ws.Cells.Item(28, 14).Value = -49326.5  # N28
$ws.Cells.Item(44, 8).Value = 950000000  # H44
$ws.Cells.Item(44, 9).Value = 950000000  # I44
$ws.Cells.Item(44, 11).Value = 950000000  # K44
$ws.Cells.Item(44, 13).Value = -949999503  # M44
$ws.Cells.Item(53, 8).Value = 59379  # H53
$ws.Cells.Item(53, 10).Value = 59379  # J53
$ws.Cells.Item(53, 12).Value = 59379  # L53
$ws.Cells.Item(53, 14).Value = -60527  # N53
$ws.Cells.Item(86, 8).Value = 41671624  # H86
$ws.Cells.Item(86, 9).Value = 83335250  # I86
$ws.Cells.Item(86, 11).Value = 83335250  # K86
$ws.Cells.Item(86, 13).Value = -83334127  # M86
$ws.Cells.Item(89, 8).Value = 41671624  # H89
$ws.Cells.Item(89, 9).Value = 83335250  # I89
$ws.Cells.Item(89, 11).Value = 416676250  # K89
$ws.Cells.Item(89, 13).Value = -416670634  # M89
$ws.Cells.Item(99, 8).Value = 9093460  # H99
$ws.Cells.Item(99, 9).Value = 2563.75  # I99
$ws.Cells.Item(99, 11).Value = 2563.75  # K99
$ws.Cells.Item(99, 13).Value = -1065.75  # M99
$ws.Cells.Item(107, 8).Value = 125004340  # H107
$ws.Cells.Item(107, 9).Value = 187503330  # I107
$ws.Cells.Item(107, 11).Value = 187503330  # K107
$ws.Cells.Item(107, 13).Value = -187501410  # M107
$ws.Cells.Item(134, 8).Value = 3314.0422  # H134
$ws.Cells.Item(134, 9).Value = 1113.3019  # I134
$ws.Cells.Item(134, 10).Value = 9794  # J134
$ws.Cells.Item(134, 11).Value = 3339.9057  # K134
$ws.Cells.Item(134, 12).Value = 29382  # L134
$ws.Cells.Item(134, 13).Value = -804.9056999999998  # M134
$ws.Cells.Item(134, 14).Value = -34452  # N134

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(6, 8).Value = 599  # H6
$ws.Cells.Item(6, 9).Value = 599  # I6
$ws.Cells.Item(6, 11).Value = 599  # K6
$ws.Cells.Item(6, 13).Value = -486  # M6
$ws.Cells.Item(16, 8).Value = 3170.6191  # H16
$ws.Cells.Item(16, 9).Value = 1108  # I16
$ws.Cells.Item(16, 11).Value = 1108  # K16
$ws.Cells.Item(16, 13).Value = -821  # M16
$ws.Cells.Item(31, 8).Value = 6021.6753  # H31
$ws.Cells.Item(31, 9).Value = 4210.387  # I31
$ws.Cells.Item(31, 10).Value = 7242.326  # J31
$ws.Cells.Item(31, 11).Value = 4210.387  # K31
$ws.Cells.Item(31, 12).Value = 7242.326  # L31
$ws.Cells.Item(31, 13).Value = -3915.387  # M31
$ws.Cells.Item(31, 14).Value = -7832.326  # N31
$ws.Cells.Item(34, 8).Value = 6021.6753  # H34
$ws.Cells.Item(34, 9).Value = 4210.387  # I34
$ws.Cells.Item(34, 10).Value = 7242.326  # J34
$ws.Cells.Item(34, 11).Value = 4210.387  # K34
$ws.Cells.Item(34, 12).Value = 7242.326  # L34
$ws.Cells.Item(34, 13).Value = -4008.387  # M34
$ws.Cells.Item(34, 14).Value = -7646.326  # N34
$ws.Cells.Item(41, 8).Value = 0  # H41
$ws.Cells.Item(41, 9).Value = 0  # I41
$ws.Cells.Item(41, 11).Value = 0  # K41
$ws.Cells.Item(41, 13).ClearContents()  # M41
$ws.Cells.Item(58, 8).Value = 10006056  # H58
$ws.Cells.Item(58, 9).Value = 25002854  # I58
$ws.Cells.Item(58, 10).Value = 8190  # J58
$ws.Cells.Item(58, 11).Value = 25002854  # K58
$ws.Cells.Item(58, 12).Value = 8190  # L58
$ws.Cells.Item(58, 13).Value = -25002651  # M58
$ws.Cells.Item(58, 14).Value = -8596  # N58
$ws.Cells.Item(76, 8).Value = 4949.3335  # H76
$ws.Cells.Item(76, 9).Value = 4949.3335  # I76
$ws.Cells.Item(76, 11).Value = 4949.3335  # K76
$ws.Cells.Item(76, 13).Value = -4634.3335  # M76
$ws.Cells.Item(79, 8).Value = 4949.3335  # H79
$ws.Cells.Item(79, 9).Value = 4949.3335  # I79
$ws.Cells.Item(79, 11).Value = 4949.3335  # K79
$ws.Cells.Item(79, 13).Value = -3857.3335  # M79
$ws.Cells.Item(107, 8).Value = 1097.5834  # H107
$ws.Cells.Item(107, 9).Value = 638.619  # I107
$ws.Cells.Item(107, 10).Value = 1740.1333  # J107
$ws.Cells.Item(107, 11).Value = 638.619  # K107
$ws.Cells.Item(107, 12).Value = 1740.1333  # L107
$ws.Cells.Item(107, 13).Value = 1281.381  # M107
$ws.Cells.Item(107, 14).Value = -5580.1333  # N107
$ws.Cells.Item(113, 8).Value = 3170.6191  # H113
$ws.Cells.Item(113, 9).Value = 1108  # I113
$ws.Cells.Item(113, 11).Value = 1108  # K113
$ws.Cells.Item(113, 13).Value = 1062  # M113
$ws.Cells.Item(132, 8).Value = 5317.343  # H132
$ws.Cells.Item(132, 9).Value = 2847.9546  # I132
$ws.Cells.Item(132, 11).Value = 8543.863799999999  # K132
$ws.Cells.Item(132, 13).Value = -6013.863799999999  # M132
$ws.Cells.Item(134, 8).Value = 5142.854  # H134
$ws.Cells.Item(134, 9).Value = 2164.0334  # I134
$ws.Cells.Item(134, 11).Value = 6492.100199999999  # K134
$ws.Cells.Item(134, 13).Value = -3957.100199999999  # M134
$ws.Cells.Item(136, 8).Value = 10006056  # H136
$ws.Cells.Item(136, 9).Value = 25002854  # I136
$ws.Cells.Item(136, 10).Value = 8190  # J136
$ws.Cells.Item(136, 11).Value = 75008562  # K136
$ws.Cells.Item(136, 12).Value = 24570  # L136
$ws.Cells.Item(136, 13).Value = -75006012  # M136
$ws.Cells.Item(136, 14).Value = -29670  # N136

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(7, 8).Value = 322.5  # H7
$ws.Cells.Item(7, 9).Value = 150  # I7
$ws.Cells.Item(7, 10).Value = 347.14285  # J7
$ws.Cells.Item(7, 11).Value = 450  # K7
$ws.Cells.Item(7, 12).Value = 1041.42855  # L7
$ws.Cells.Item(7, 13).Value = -338  # M7
$ws.Cells.Item(7, 14).Value = -1265.42855  # N7
$ws.Cells.Item(12, 8).Value = 2941952.2  # H12
$ws.Cells.Item(12, 9).Value = 668.375  # I12
$ws.Cells.Item(12, 10).Value = 5556427  # J12
$ws.Cells.Item(12, 11).Value = 2005.125  # K12
$ws.Cells.Item(12, 12).Value = 16669281  # L12
$ws.Cells.Item(12, 13).Value = -1832.125  # M12
$ws.Cells.Item(12, 14).Value = -16669627  # N12
$ws.Cells.Item(56, 8).Value = 6098.8  # H56
$ws.Cells.Item(56, 9).Value = 6098.8  # I56
$ws.Cells.Item(56, 11).Value = 6098.8  # K56
$ws.Cells.Item(56, 13).Value = -5568.8  # M56
$ws.Cells.Item(64, 8).Value = 1783.8  # H64
$ws.Cells.Item(64, 9).Value = 709.5  # I64
$ws.Cells.Item(64, 11).Value = 2128.5  # K64
$ws.Cells.Item(64, 13).Value = -1858.5  # M64
$ws.Cells.Item(67, 8).Value = 1783.8  # H67
$ws.Cells.Item(67, 9).Value = 709.5  # I67
$ws.Cells.Item(67, 11).Value = 2128.5  # K67
$ws.Cells.Item(67, 13).Value = -1192.5  # M67
$ws.Cells.Item(68, 8).Value = 2873.2263  # H68
$ws.Cells.Item(68, 10).Value = 2925.6956  # J68
$ws.Cells.Item(68, 12).Value = 8777.086800000001  # L68
$ws.Cells.Item(68, 14).Value = -10399.0868  # N68
$ws.Cells.Item(71, 8).Value = 2873.2263  # H71
$ws.Cells.Item(71, 10).Value = 2925.6956  # J71
$ws.Cells.Item(71, 12).Value = 26331.2604  # L71
$ws.Cells.Item(71, 14).Value = -34443.2604  # N71
$ws.Cells.Item(93, 8).Value = 6771.909  # H93
$ws.Cells.Item(93, 9).Value = 2294  # I93
$ws.Cells.Item(93, 11).Value = 6882  # K93
$ws.Cells.Item(93, 13).Value = -5010  # M93
$ws.Cells.Item(107, 8).Value = 14286336  # H107
$ws.Cells.Item(107, 10).Value = 16667325  # J107
$ws.Cells.Item(107, 12).Value = 50001975  # L107
$ws.Cells.Item(107, 14).Value = -50005815  # N107
$ws.Cells.Item(113, 8).Value = 1970.0238  # H113
$ws.Cells.Item(113, 9).Value = 1097.8182  # I113
$ws.Cells.Item(113, 11).Value = 3293.4546  # K113
$ws.Cells.Item(113, 13).Value = -1123.4546  # M113
$ws.Cells.Item(117, 8).Value = 625770  # H117
$ws.Cells.Item(117, 9).Value = 29  # I117
$ws.Cells.Item(117, 10).Value = 715161.5600000001  # J117
$ws.Cells.Item(117, 11).Value = 87  # K117
$ws.Cells.Item(117, 12).Value = 2145484.68  # L117
$ws.Cells.Item(117, 13).Value = 3355  # M117
$ws.Cells.Item(117, 14).Value = -2152368.68  # N117
$ws.Cells.Item(119, 8).Value = 10000  # H119
$ws.Cells.Item(119, 9).Value = 10000  # I119
$ws.Cells.Item(119, 11).Value = 30000  # K119
$ws.Cells.Item(119, 13).Value = -25162  # M119
$ws.Cells.Item(131, 8).Value = 1734.875  # H131
$ws.Cells.Item(131, 10).Value = 2040.9656  # J131
$ws.Cells.Item(131, 12).Value = 6122.8968  # L131
$ws.Cells.Item(131, 14).Value = -16202.8968  # N131
$ws.Cells.Item(137, 8).Value = 176999.67  # H137
$ws.Cells.Item(137, 10).Value = 102181.55  # J137
$ws.Cells.Item(137, 12).Value = 306544.65  # L137
$ws.Cells.Item(137, 14).Value = -316744.65  # N137

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(52, 8).Value = 90000  # H52
$ws.Cells.Item(52, 10).Value = 90000  # J52
$ws.Cells.Item(52, 12).Value = 90000  # L52
$ws.Cells.Item(52, 14).Value = -90518  # N52
$ws.Cells.Item(80, 8).Value = 4754.5557  # H80
$ws.Cells.Item(80, 9).Value = 4673.875  # I80
$ws.Cells.Item(80, 10).Value = 5400  # J80
$ws.Cells.Item(80, 11).Value = 4673.875  # K80
$ws.Cells.Item(80, 12).Value = 5400  # L80
$ws.Cells.Item(80, 13).Value = -3675.875  # M80
$ws.Cells.Item(80, 14).Value = -7396  # N80
$ws.Cells.Item(83, 8).Value = 4754.5557  # H83
$ws.Cells.Item(83, 9).Value = 4673.875  # I83
$ws.Cells.Item(83, 10).Value = 5400  # J83
$ws.Cells.Item(83, 11).Value = 23369.375  # K83
$ws.Cells.Item(83, 12).Value = 27000  # L83
$ws.Cells.Item(83, 13).Value = -18377.375  # M83
$ws.Cells.Item(83, 14).Value = -36984  # N83
$ws.Cells.Item(113, 8).Value = 3012.0303  # H113
$ws.Cells.Item(113, 9).Value = 1964.9166  # I113
$ws.Cells.Item(113, 11).Value = 1964.9166  # K113
$ws.Cells.Item(113, 13).Value = 205.0834  # M113
$ws.Cells.Item(122, 8).Value = 89349.336  # H122
$ws.Cells.Item(122, 9).Value = 170867.33  # I122
$ws.Cells.Item(122, 11).Value = 512601.99  # K122
$ws.Cells.Item(122, 13).Value = -510151.99  # M122
$ws.Cells.Item(132, 8).Value = 5503.952  # H132
$ws.Cells.Item(132, 10).Value = 9510  # J132
$ws.Cells.Item(132, 12).Value = 28530  # L132
$ws.Cells.Item(132, 14).Value = -33590  # N132

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(13, 8).Value = 0  # H13
$ws.Cells.Item(13, 9).Value = 0  # I13
$ws.Cells.Item(13, 11).Value = 0  # K13
$ws.Cells.Item(13, 13).ClearContents()  # M13
$ws.Cells.Item(35, 8).Value = 24628  # H35
$ws.Cells.Item(35, 9).Value = 24628  # I35
$ws.Cells.Item(35, 11).Value = 24628  # K35
$ws.Cells.Item(35, 13).Value = -24292  # M35
$ws.Cells.Item(40, 8).Value = 6731.625  # H40
$ws.Cells.Item(40, 9).Value = 6445.1816  # I40
$ws.Cells.Item(40, 10).Value = 7361.8  # J40
$ws.Cells.Item(40, 11).Value = 6445.1816  # K40
$ws.Cells.Item(40, 12).Value = 7361.8  # L40
$ws.Cells.Item(40, 13).Value = -6309.1816  # M40
$ws.Cells.Item(40, 14).Value = -7633.8  # N40
$ws.Cells.Item(43, 8).Value = 7000  # H43
$ws.Cells.Item(43, 9).Value = 7000  # I43
$ws.Cells.Item(43, 11).Value = 7000  # K43
$ws.Cells.Item(43, 13).Value = -6807  # M43
$ws.Cells.Item(82, 8).Value = 504826.94  # H82
$ws.Cells.Item(82, 9).Value = 940285.25  # I82
$ws.Cells.Item(82, 10).Value = 2375  # J82
$ws.Cells.Item(82, 11).Value = 940285.25  # K82
$ws.Cells.Item(82, 12).Value = 2375  # L82
$ws.Cells.Item(82, 13).Value = -939924.25  # M82
$ws.Cells.Item(82, 14).Value = -3097  # N82
$ws.Cells.Item(85, 8).Value = 504826.94  # H85
$ws.Cells.Item(85, 9).Value = 940285.25  # I85
$ws.Cells.Item(85, 10).Value = 2375  # J85
$ws.Cells.Item(85, 11).Value = 940285.25  # K85
$ws.Cells.Item(85, 12).Value = 2375  # L85
$ws.Cells.Item(85, 13).Value = -939037.25  # M85
$ws.Cells.Item(85, 14).Value = -4871  # N85
$ws.Cells.Item(122, 8).Value = 6736.4517  # H122
$ws.Cells.Item(122, 9).Value = 6343.75  # I122
$ws.Cells.Item(122, 10).Value = 7450.4546  # J122
$ws.Cells.Item(122, 11).Value = 19031.25  # K122
$ws.Cells.Item(122, 12).Value = 22351.3638  # L122
$ws.Cells.Item(122, 13).Value = -16581.25  # M122
$ws.Cells.Item(122, 14).Value = -27251.3638  # N122

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(15, 8).Value = 25000  # H15
$ws.Cells.Item(15, 9).Value = 25000  # I15
$ws.Cells.Item(15, 11).Value = 25000  # K15
$ws.Cells.Item(15, 13).Value = -24712  # M15
$ws.Cells.Item(81, 8).Value = 26261612  # H81
$ws.Cells.Item(81, 9).Value = 3334665.8  # I81
$ws.Cells.Item(81, 10).Value = 40017780  # J81
$ws.Cells.Item(81, 11).Value = 6669331.6  # K81
$ws.Cells.Item(81, 12).Value = 80035560  # L81
$ws.Cells.Item(81, 13).Value = -6668270.6  # M81
$ws.Cells.Item(81, 14).Value = -80037682  # N81
$ws.Cells.Item(84, 8).Value = 26261612  # H84
$ws.Cells.Item(84, 9).Value = 3334665.8  # I84
$ws.Cells.Item(84, 10).Value = 40017780  # J84
$ws.Cells.Item(84, 11).Value = 33346658  # K84
$ws.Cells.Item(84, 12).Value = 400177800  # L84
$ws.Cells.Item(84, 13).Value = -33341354  # M84
$ws.Cells.Item(84, 14).Value = -400188408  # N84
$ws.Cells.Item(107, 8).Value = 767.75  # H107
$ws.Cells.Item(107, 9).Value = 433.66666  # I107
$ws.Cells.Item(107, 11).Value = 1300.99998  # K107
$ws.Cells.Item(107, 13).Value = 619.0000199999999  # M107
$ws.Cells.Item(122, 8).Value = 3637.303  # H122
$ws.Cells.Item(122, 9).Value = 3001.0688  # I122
$ws.Cells.Item(122, 11).Value = 9003.206399999999  # K122
$ws.Cells.Item(122, 13).Value = -6553.206399999999  # M122
$ws.Cells.Item(132, 8).Value = 23434.62  # H132
$ws.Cells.Item(132, 9).Value = 9527.069  # I132
$ws.Cells.Item(132, 11).Value = 28581.207  # K132
$ws.Cells.Item(132, 13).Value = -26051.207  # M132
$ws.Cells.Item(135, 8).Value = 0  # H135
$ws.Cells.Item(135, 10).Value = 0  # J135
$ws.Cells.Item(135, 12).ClearContents()  # L135
$ws.Cells.Item(135, 14).ClearContents()  # N135
$ws.Cells.Item(136, 8).Value = 43482320  # H136
$ws.Cells.Item(136, 9).Value = 62502584  # I136
$ws.Cells.Item(136, 11).Value = 187507752  # K136
$ws.Cells.Item(136, 13).Value = -187505202  # M136
$ws.Cells.Item(137, 8).Value = 0  # H137
$ws.Cells.Item(137, 10).Value = 0  # J137
$ws.Cells.Item(137, 12).ClearContents()  # L137
$ws.Cells.Item(137, 14).ClearContents()  # N137
